$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update component forecast values (rows 4-73, column B)
$ws.Range("B4").Value = 1.960591199849219
$ws.Range("B5").Value = -0.2410745035013093
$ws.Range("B6").Value = -2.989151587480549
$ws.Range("B7").Value = -1.992436460080896
$ws.Range("B8").Value = -2.911016840458359
$ws.Range("B9").Value = 0.5161858965418986
$ws.Range("B10").Value = 3.496635313038965
$ws.Range("B11").Value = 2.022610195364777
$ws.Range("B12").Value = 4.001048740912381
$ws.Range("B13").Value = 2.471556584906722
$ws.Range("B14").Value = 1.809198400323567
$ws.Range("B15").Value = 0.9930361476913703
$ws.Range("B16").Value = 2.297598938510475
$ws.Range("B17").Value = 1.683866211357127
$ws.Range("B18").Value = -0.2000000000000028
$ws.Range("B19").Value = 0
$ws.Range("B20").Value = 0.4230141822949065
$ws.Range("B21").Value = 0.1000000000000085
$ws.Range("B22").Value = -0.8135079514120491
$ws.Range("B23").Value = 0.5011757240633443
$ws.Range("B24").Value = 1.792421234023649
$ws.Range("B25").Value = 1.000583728414611
$ws.Range("B26").Value = 1.700400414966595
$ws.Range("B27").Value = 1.466172678153498
$ws.Range("B28").Value = 1.081988526104823
$ws.Range("B29").Value = 1.353592454001856
$ws.Range("B30").Value = 0.4035607671893757
$ws.Range("B31").Value = 0.9949231225214561
$ws.Range("B32").Value = 1.476517947415175
$ws.Range("B33").Value = 1.292287545731185
$ws.Range("B34").Value = 0.4919486237229052
$ws.Range("B35").Value = 1.268067508957756
$ws.Range("B36").Value = 0.9496767498412169
$ws.Range("B37").Value = 1.183775223599739
$ws.Range("B38").Value = 1.014999674053939
$ws.Range("B39").Value = 1.191228972472018
$ws.Range("B40").Value = 1.836987253553829
$ws.Range("B41").Value = 1.33483377666461
$ws.Range("B42").Value = 1.622456305393911
$ws.Range("B43").Value = 1.303119407694879
$ws.Range("B44").Value = 0.3932797567441497
$ws.Range("B45").Value = 1.381065482681691
$ws.Range("B46").Value = 0.7
$ws.Range("B47").Value = 1.712723949918967
$ws.Range("B48").Value = -0.9
$ws.Range("B49").Value = 0.914789154762218
$ws.Range("B50").Value = 1.152379855430368
$ws.Range("B51").Value = -1.558195694240041
$ws.Range("B52").Value = -24.4
$ws.Range("B53").Value = 12.42550598425463
$ws.Range("B54").Value = 2.813188392915293
$ws.Range("B55").Value = 1.806909003397877
$ws.Range("B56").Value = 1.087285196410733
$ws.Range("B57").Value = 1.450779288666709
$ws.Range("B58").Value = 0.368631518524424
$ws.Range("B59").Value = -2.021680416268424
$ws.Range("B60").Value = 1.118165634023697
$ws.Range("B61").Value = 0.1289585187160185
$ws.Range("B62").Value = -0.4259770459179748
$ws.Range("B63").Value = -0.6379881427730965
$ws.Range("B64").Value = -0.04993664583679447
$ws.Range("B65").Value = 0.3662234592800075
$ws.Range("B66").Value = -1.023919595764212
$ws.Range("B67").Value = -0.2547854660834332
$ws.Range("B68").Value = 1.565950786385088
$ws.Range("B69").Value = -0.2475555093771362
$ws.Range("B70").Value = 0.386705041387188
$ws.Range("B71").Value = -0.6950112442994083
$ws.Range("B72").Value = -0.3135696110003181
$ws.Range("B73").Value = -0.02110736963342674

# Remove the trailing rows (74-82) that are no longer part of the series
$ws.Range("A74:B82").EntireRow.Delete()
